# Add a new test case row (row 31) to the "protocoltestcasedetails" sheet:
# Sno.=30, test_case_name=testcase30_csv_csv_3mill50cols_content,
# test_case_file_path=CONCAT formula, protocol_application_name(D)="N"

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("protocoltestcasedetails")

$ws.Range("A31").Value = 30
$ws.Range("B31").Value = "testcase30_csv_csv_3mill50cols_content"
$ws.Range("C31").Formula = '=_xlfn.CONCAT("/app/test/testcases/",B31,".xlsx")'
$ws.Range("D31").Value = "N"
